# Refresh of the live crypto-price snapshot (coinranking.com feed).
# Updates Price (D) / Volume(1h) (E) per coin, plus two row swaps where
# the source feed re-ordered neighbouring coins (Aptos/Algorand,
# Decentraland/NEARProtocol) on rows 39-40 and 45-46.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.265.72'
$ws.Range("E2").Value = '  +0.13%  '
# Row 3
$ws.Range("D3").Value = '2.002.24'
$ws.Range("E3").Value = '  +5.65%  '
# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.05%  '
# Row 5
$ws.Range("D5").Value = '''323.75'
$ws.Range("E5").Value = '  +0.78%  '
# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.06%  '
# Row 7
$ws.Range("D7").Value = '''0.5103'
$ws.Range("E7").Value = '  +0.78%  '
# Row 8
$ws.Range("D8").Value = '''0.4138'
$ws.Range("E8").Value = '  +2.33%  '
# Row 9
$ws.Range("D9").Value = '''0.08688'
$ws.Range("E9").Value = '  +4.71%  '
# Row 10
$ws.Range("D10").Value = '''1.137'
$ws.Range("E10").Value = '  +2.34%  '
# Row 11
$ws.Range("D11").Value = '''42.80'
$ws.Range("E11").Value = '  +1.36%  '
# Row 12
$ws.Range("D12").Value = '''24.98'
$ws.Range("E12").Value = '  +3.19%  '
# Row 13
$ws.Range("D13").Value = '2.000.75'
$ws.Range("E13").Value = '  +7.03%  '
# Row 14
$ws.Range("D14").Value = '''6.533'
$ws.Range("E14").Value = '  +1.80%  '
# Row 15
$ws.Range("D15").Value = '''7.431'
$ws.Range("E15").Value = '  +1.44%  '
# Row 16
$ws.Range("E16").Value = '  +0.04%  '
# Row 17
$ws.Range("D17").Value = '''94.17'
$ws.Range("E17").Value = '  +1.44%  '
# Row 18
$ws.Range("D18").Value = '''0.00001117'
# Row 19
$ws.Range("D19").Value = '''0.06546'
$ws.Range("E19").Value = '  +1.00%  '
# Row 20
$ws.Range("D20").Value = '''18.93'
$ws.Range("E20").Value = '  +2.76%  '
# Row 21
$ws.Range("D21").Value = '''1.000'
$ws.Range("E21").Value = '  +0.06%  '
# Row 22
$ws.Range("D22").Value = '''6.144'
$ws.Range("E22").Value = '  +3.72%  '
# Row 23
$ws.Range("D23").Value = '30.322.41'
$ws.Range("E23").Value = '  +0.30%  '
# Row 24
$ws.Range("D24").Value = '''11.70'
$ws.Range("E24").Value = '  +3.53%  '
# Row 25
$ws.Range("D25").Value = '''2.210'
$ws.Range("E25").Value = '  +1.13%  '
# Row 26
$ws.Range("D26").Value = '2.233.86'
$ws.Range("E26").Value = '  +6.84%  '
# Row 27
$ws.Range("D27").Value = '''22.57'
$ws.Range("E27").Value = '  +4.36%  '
# Row 28
$ws.Range("D28").Value = '''163.53'
$ws.Range("E28").Value = '  +1.64%  '
# Row 29
$ws.Range("D29").Value = '''2.398'
$ws.Range("E29").Value = '  +5.51%  '
# Row 30
$ws.Range("D30").Value = '''131.36'
$ws.Range("E30").Value = '  +1.90%  '
# Row 31
$ws.Range("D31").Value = '''1.141'
$ws.Range("E31").Value = '  +3.24%  '
# Row 32
$ws.Range("D32").Value = '''0.1052'
$ws.Range("E32").Value = '  +0.95%  '
# Row 33
$ws.Range("E33").Value = '  +1.02%  '
# Row 34
$ws.Range("D34").Value = '''3.835'
$ws.Range("E34").Value = '  +3.41%  '
# Row 35
$ws.Range("D35").Value = '''1.340'
$ws.Range("E35").Value = '  +12.70%  '
# Row 36
$ws.Range("D36").Value = '''0.02499'
$ws.Range("E36").Value = '  +1.91%  '
# Row 37
$ws.Range("E37").Value = '  +1.40%  '
# Row 38
$ws.Range("D38").Value = '''0.06581'
$ws.Range("E38").Value = '  +1.96%  '
# Row 39
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '''12.35'
$ws.Range("E39").Value = '  +8.23%  '
# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '''0.2201'
$ws.Range("E40").Value = '  +2.03%  '
# Row 41
$ws.Range("D41").Value = '''8.914'
$ws.Range("E41").Value = '  +3.60%  '
# Row 42
$ws.Range("D42").Value = '''0.6609'
$ws.Range("E42").Value = '  +3.11%  '
# Row 43
$ws.Range("E43").Value = '  +0.86%  '
# Row 44
$ws.Range("D44").Value = '''13.61'
$ws.Range("E44").Value = '  +3.11%  '
# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.6151'
$ws.Range("E45").Value = '  +2.93%  '
# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''2.209'
$ws.Range("E46").Value = '  +3.10%  '
# Row 47
$ws.Range("D47").Value = '''3.660'
$ws.Range("E47").Value = '  +0.51%  '
# Row 48
$ws.Range("D48").Value = '''1.263'
$ws.Range("E48").Value = '  +3.98%  '
# Row 49
$ws.Range("D49").Value = '''124.59'
$ws.Range("E49").Value = '  +0.77%  '
# Row 50
$ws.Range("D50").Value = '''79.95'
$ws.Range("E50").Value = '  +1.43%  '
# Row 51
$ws.Range("D51").Value = '''0.06890'
$ws.Range("E51").Value = '  +1.80%  '
